$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D/E columns so numeric-looking strings
# (e.g. "0.9993", "1.103") are not coerced into real numbers by COM,
# matching the source workbook where these are plain text cells.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "D46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.724.92"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.758.39"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "324.95"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4610"
$ws.Range("E7").Value = "  +7.88%  "
$ws.Range("D8").Value = "0.3609"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "0.07539"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "42.22"
$ws.Range("D11").Value = "1.103"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "0.9988"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "20.90"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "6.029"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "7.133"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "1.755.32"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "92.82"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "0.06429"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "16.84"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "5.832"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").Value = "27.754.89"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "11.26"
$ws.Range("D25").Value = "2.108"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "163.65"
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("D27").Value = "20.52"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "1.952.91"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "2.104"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "127.19"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").Value = "1.080"
$ws.Range("E31").Value = "  -6.87%  "
$ws.Range("D32").Value = "0.09217"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D33").Value = "3.665"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "5.560"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").Value = "11.98"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").Value = "0.02305"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "0.2105"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "0.06058"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "0.6398"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "4.986"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "1.204"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("D42").Value = "1.381"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "7.841"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "13.33"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "0.5930"
$ws.Range("D46").Value = "3.716"
$ws.Range("D47").Value = "123.28"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "1.965"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").Value = "1.150"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "0.06880"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "72.46"
$ws.Range("E51").Value = "  -2.17%  "
